# Apply roster corrections to the "Golden State" sheet:
#  1. "Anthony Lamb (TW)" -> "Anthony Lamb" (two-way tag dropped)
#  2. Row 11 and Row 12 swap their player data (Stephen Curry moves up to
#     row 11 / jersey #30, Ty Jerome moves down to row 12 / jersey #10),
#     while the "No." rank column (A) stays as-is.
#  3. "Lester Quinones" -> "Lester Quinones (TW)" (two-way tag added)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Anthony Lamb loses the "(TW)" suffix
$ws.Range("C9").Value = "Anthony Lamb"

# 2. Swap the full row-11 / row-12 player records (No. column A untouched)
$ws.Range("B11").Value = 30
$ws.Range("C11").Value = "Stephen Curry"
$ws.Range("D11").Value = "PG"
$ws.Range("E11").Value = "6-2"
$ws.Range("F11").Value = 185
$ws.Range("G11").Value = "March 14, 1988"
$ws.Range("I11").Value = "13"
$ws.Range("J11").Value = "Davidson"
$ws.Range("K11").Value = "https://www.basketball-reference.com/players/c/curryst01.html"

$ws.Range("B12").Value = 10
$ws.Range("C12").Value = "Ty Jerome (TW)"
$ws.Range("D12").Value = "SG"
$ws.Range("E12").Value = "6-5"
$ws.Range("F12").Value = 195
$ws.Range("G12").Value = "July 8, 1997"
$ws.Range("I12").Value = "3"
$ws.Range("J12").Value = "Virginia"
$ws.Range("K12").Value = "https://www.basketball-reference.com/players/j/jeromty01.html"

# 3. Lester Quinones gains the "(TW)" suffix
$ws.Range("C17").Value = "Lester Quinones (TW)"
